$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells whose new value is plain numeric text so Excel
# keeps them as text (matching the source inline-string cells) instead of
# auto-converting to a number when the value is assigned below.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('D2').Value = '27.147.32'
$ws.Range('E2').Value = '  -1.87%  '

$ws.Range('D3').Value = '1.559.52'
$ws.Range('E3').Value = '  -1.78%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '206.62'
$ws.Range('E5').Value = '  -0.29%  '

$ws.Range('E6').Value = '  -1.08%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = '22.39'
$ws.Range('E8').Value = '  +0.78%  '

$ws.Range('E9').Value = '  -1.87%  '

$ws.Range('E10').Value = '  +0.20%  '

$ws.Range('D11').Value = '0.0861'
$ws.Range('E11').Value = '  -0.71%  '

$ws.Range('D12').Value = '1.782.05'
$ws.Range('E12').Value = '  -1.71%  '

$ws.Range('D13').Value = '1.561.27'
$ws.Range('E13').Value = '  -1.75%  '

$ws.Range('E14').Value = '  -2.18%  '

$ws.Range('D15').Value = '0.515'
$ws.Range('E15').Value = '  -2.82%  '

$ws.Range('D16').Value = '62.83'
$ws.Range('E16').Value = '  -0.95%  '

$ws.Range('D17').Value = '27.156.63'

$ws.Range('D18').Value = '213.32'
$ws.Range('E18').Value = '  -2.78%  '

$ws.Range('D19').Value = '0.0₃0687'
$ws.Range('E19').Value = '  -1.13%  '

$ws.Range('E20').Value = '  -1.19%  '

$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  -0.57%  '

$ws.Range('D23').Value = '9.35'
$ws.Range('E23').Value = '  -3.11%  '

$ws.Range('E24').Value = '  -0.09%  '

$ws.Range('D25').Value = '152.11'
$ws.Range('E25').Value = '  -0.93%  '

$ws.Range('D26').Value = '6.60'
$ws.Range('E26').Value = '  -3.26%  '

$ws.Range('D27').Value = '14.89'
$ws.Range('E27').Value = '  -1.50%  '

$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.01%  '

$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '0.104'
$ws.Range('E29').Value = '  -1.13%  '

$ws.Range('E30').Value = '  -0.52%  '

$ws.Range('D31').Value = '0.0462'
$ws.Range('E31').Value = '  -1.14%  '

$ws.Range('E32').Value = '  -1.78%  '

$ws.Range('D33').Value = '1.384.57'
$ws.Range('E33').Value = '  +1.16%  '

$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  +0.36%  '

$ws.Range('E35').Value = '  +0.36%  '

$ws.Range('D36').Value = '0.946'
$ws.Range('E36').Value = '  -3.03%  '

$ws.Range('E37').Value = '  -1.32%  '

$ws.Range('E38').Value = '  -1.11%  '

$ws.Range('D39').Value = '0.814'
$ws.Range('E39').Value = '  -1.31%  '

$ws.Range('D40').Value = '0.517'
$ws.Range('E40').Value = '  -3.26%  '

$ws.Range('E41').Value = '  +0.03%  '

$ws.Range('E42').Value = '  +1.56%  '

$ws.Range('D43').Value = '1.78'
$ws.Range('E43').Value = '  +4.05%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '2.17'
$ws.Range('E44').Value = '  +0.17%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '63.33'
$ws.Range('E45').Value = '  -1.25%  '

$ws.Range('E46').Value = '  +1.06%  '

$ws.Range('D47').Value = '1.695.11'
$ws.Range('E47').Value = '  -1.68%  '

$ws.Range('D48').Value = '85.57'
$ws.Range('E48').Value = '  -2.32%  '

$ws.Range('D49').Value = '0.0₇0992'
$ws.Range('E49').Value = '  -1.42%  '

$ws.Range('D50').Value = '0.0493'
$ws.Range('E50').Value = '  -0.46%  '

$ws.Range('E51').Value = '  +0.07%  '
